$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Replace the "LmN" text labels in column B with plain integers 1..20
for ($i = 1; $i -le 20; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $i
}

# Reflect the selection that was active when the workbook was saved
$ws.Range("B22").Select()
